# Updates US copy to commit #5e209809
$wb = $excel.ActiveWorkbook
$wsAbout  = $wb.Worksheets.Item("About")
$wsCRpUNL = $wb.Worksheets.Item("CRpUNL")

# ---------------------------------------------------------------------------
# CRpUNL sheet: header relabeled from a fraction to an absolute MW value,
# and the unit note updated accordingly.
# ---------------------------------------------------------------------------
$wsCRpUNL.Range("B1").Value = "MW retired"
$wsCRpUNL.Range("A1").Value = "Unit: MW/($/MW)"

# ---------------------------------------------------------------------------
# About sheet: update the "These includes..." note text, and append two new
# explanatory lines about biomass plants (rows 16-17).
# ---------------------------------------------------------------------------
$wsAbout.Range("A10").Value = "These includes: natural gas steam turbines and petroleum plants. For these plant types we set the "
$wsAbout.Range("A16").Value = "Likewise, biomass plants are often colocated with cheap supply and part of integrated"
$wsAbout.Range("A17").Value = "CHP or industrial systems, and we therefore do not subject them to economic retirement."

# Numeric value updates (column B) -----------------------------------------
$wsCRpUNL.Range("B2").Value  = 0.03    # hard coal
$wsCRpUNL.Range("B3").Value  = 0.03    # natural gas steam turbine
$wsCRpUNL.Range("B4").Value  = 0.03    # natural gas combined cycle
$wsCRpUNL.Range("B5").Value  = 0.03    # nuclear
$wsCRpUNL.Range("B7").Value  = 0.03    # onshore wind
$wsCRpUNL.Range("B8").Value  = 0.03    # solar PV
$wsCRpUNL.Range("B9").Value  = 0       # solar thermal
$wsCRpUNL.Range("B10").Value = 0       # biomass
$wsCRpUNL.Range("B11").Value = 0       # geothermal
$wsCRpUNL.Range("B12").Value = 0       # petroleum
$wsCRpUNL.Range("B13").Value = 0.03    # natural gas peaker
$wsCRpUNL.Range("B14").Value = 0.03    # lignite
$wsCRpUNL.Range("B15").Value = 0.03    # offshore wind
$wsCRpUNL.Range("B18").Value = 0       # municipal solid waste
$wsCRpUNL.Range("B19").Value = 0.03    # hard coal w CCS
$wsCRpUNL.Range("B20").Value = 0.03    # natural gas combined cycle w CCS
$wsCRpUNL.Range("B21").Value = 0.03    # biomass w CCS
$wsCRpUNL.Range("B22").Value = 0.03    # lignite w CCS
$wsCRpUNL.Range("B23").Value = 0.03    # small modular reactor
$wsCRpUNL.Range("B24").Value = 0.03    # hydrogen combustion turbine
$wsCRpUNL.Range("B25").Value = 0.03    # hydrogen combined cycle

# Cells that are now exactly zero get the shaded "not applicable" fill style
# that already exists on B6/B16/B17 (fillId referencing the grey theme fill).
$wsCRpUNL.Range("B6").Copy()
$wsCRpUNL.Range("B9").PasteSpecial(-4122)
$wsCRpUNL.Range("B10").PasteSpecial(-4122)
$wsCRpUNL.Range("B11").PasteSpecial(-4122)
$wsCRpUNL.Range("B12").PasteSpecial(-4122)
$wsCRpUNL.Range("B18").PasteSpecial(-4122)

[void]$wsCRpUNL.Range("D14").Select()

# Restore "About" as the active/selected sheet & cell, matching the original
# workbook state (only the selected cell moved from A15 to A18).
[void]$wsAbout.Activate()
[void]$wsAbout.Range("A18").Select()
